# Insert two new weekly price rows (Coliflor, Terminal Hortofrutícola Agro Chillán)
# at row 456, pushing the existing rows 456-474 down to 458-476.
# The new rows correspond to date 2023-05-29 (serial 45075), one "Primera" and one
# "Segunda" quality record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 456 onward down by inserting two blank rows at 456:457.
$ws.Range("456:457").Insert()

# --- Row 456: new "Primera" quality record ---
$ws.Cells.Item(456, 1).Value = 7
$ws.Cells.Item(456, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(456, 3).Value = "Ñuble"
$ws.Cells.Item(456, 4).Value = 45075
$ws.Cells.Item(456, 5).Value = 16
$ws.Cells.Item(456, 6).Value = 100112008
$ws.Cells.Item(456, 7).Value = "Coliflor"
$ws.Cells.Item(456, 8).Value = "Sin especificar"
$ws.Cells.Item(456, 9).Value = "Primera"
$ws.Cells.Item(456, 10).Value = 150
$ws.Cells.Item(456, 11).Value = 1200
$ws.Cells.Item(456, 12).Value = 1200
$ws.Cells.Item(456, 13).Value = 1200
$ws.Cells.Item(456, 14).Value = "$/unidad"
$ws.Cells.Item(456, 15).Value = "Región del Maule"
$ws.Cells.Item(456, 16).Value = 1200
$ws.Cells.Item(456, 17).Value = 1
$ws.Cells.Item(456, 18).Value = "Hortaliza"

# --- Row 457: new "Segunda" quality record ---
$ws.Cells.Item(457, 1).Value = 7
$ws.Cells.Item(457, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(457, 3).Value = "Ñuble"
$ws.Cells.Item(457, 4).Value = 45075
$ws.Cells.Item(457, 5).Value = 16
$ws.Cells.Item(457, 6).Value = 100112008
$ws.Cells.Item(457, 7).Value = "Coliflor"
$ws.Cells.Item(457, 8).Value = "Sin especificar"
$ws.Cells.Item(457, 9).Value = "Segunda"
$ws.Cells.Item(457, 10).Value = 150
$ws.Cells.Item(457, 11).Value = 1000
$ws.Cells.Item(457, 12).Value = 1000
$ws.Cells.Item(457, 13).Value = 1000
$ws.Cells.Item(457, 14).Value = "$/unidad"
$ws.Cells.Item(457, 15).Value = "Región del Maule"
$ws.Cells.Item(457, 16).Value = 1000
$ws.Cells.Item(457, 17).Value = 1
$ws.Cells.Item(457, 18).Value = "Hortaliza"
